# Fix mangled UTF-8 "±" (plus-minus) that was double-encoded as "Â±"
# (UTF-8 bytes of U+00B1 misread as Latin-1/CP1252) back to the correct
# "±" character. Affects the metric cells (B2:H17) on the f1_score,
# training_time and test_time sheets.

$wb = $excel.ActiveWorkbook
$sheetNames = @("f1_score", "training_time", "test_time")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    for ($r = 2; $r -le 17; $r++) {
        for ($c = 2; $c -le 8; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            $v = $cell.Value()
            if ($v -ne $null) {
                $fixed = $v.Replace("Â±", "±")
                if ($fixed -ne $v) {
                    $cell.Value = $fixed
                }
            }
        }
    }
}
